# Auto-generated edit script: applies scheduled market-data refresh values
# to the profit-calculation columns (H-N) across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 3308.926
$ws.Range("I28").Value = 2440.9524
$ws.Range("J28").Value = 6346.8335
$ws.Range("K28").Value = 2440.9524
$ws.Range("L28").Value = 6346.8335
$ws.Range("M28").Value = -1955.9524
$ws.Range("N28").Value = -7316.8335
$ws.Range("H43").Value = 2362.6667
$ws.Range("I43").Value = 1299
$ws.Range("J43").Value = 2495.625
$ws.Range("K43").Value = 1299
$ws.Range("L43").Value = 2495.625
$ws.Range("M43").Value = -1230
$ws.Range("N43").Value = -2633.625
$ws.Range("H70").Value = 4469.4595
$ws.Range("I70").Value = 1325.8636
$ws.Range("J70").Value = 9080.066000000001
$ws.Range("K70").Value = 3977.5908
$ws.Range("L70").Value = 27240.198
$ws.Range("M70").Value = -3707.5908
$ws.Range("N70").Value = -27780.198
$ws.Range("H73").Value = 4469.4595
$ws.Range("I73").Value = 1325.8636
$ws.Range("J73").Value = 9080.066000000001
$ws.Range("K73").Value = 3977.5908
$ws.Range("L73").Value = 27240.198
$ws.Range("M73").Value = -3041.5908
$ws.Range("N73").Value = -29112.198
$ws.Range("H112").Value = 2161.7144
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 2342.182
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 7026.545999999999
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -9242.545999999998
$ws.Range("H137").Value = 2014.4
$ws.Range("I137").Value = 1383.7059
$ws.Range("J137").Value = 2610.0557
$ws.Range("K137").Value = 4151.1177
$ws.Range("L137").Value = 7830.1671
$ws.Range("M137").Value = -1601.1177
$ws.Range("N137").Value = -12930.1671
$ws.Range("H138").Value = 3267.2676
$ws.Range("I138").Value = 2707.4194
$ws.Range("J138").Value = 3701.15
$ws.Range("K138").Value = 8122.2582
$ws.Range("L138").Value = 11103.45
$ws.Range("M138").Value = -2982.2582
$ws.Range("N138").Value = -21383.45
$ws.Range("H141").Value = 4381.3794
$ws.Range("I141").Value = 3466.4644
$ws.Range("J141").Value = 29999
$ws.Range("K141").Value = 10399.3932
$ws.Range("L141").Value = 89997
$ws.Range("M141").Value = -5219.393199999999
$ws.Range("N141").Value = -100357

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2334.4443
$ws.Range("I61").Value = 2708.7334
$ws.Range("J61").Value = 463
$ws.Range("K61").Value = 2708.7334
$ws.Range("L61").Value = 463
$ws.Range("M61").Value = -2496.7334
$ws.Range("N61").Value = -887
$ws.Range("H74").Value = 2161.925
$ws.Range("I74").Value = 1437.3077
$ws.Range("K74").Value = 1437.3077
$ws.Range("M74").Value = -563.3077000000001
$ws.Range("H77").Value = 2161.925
$ws.Range("I77").Value = 1437.3077
$ws.Range("K77").Value = 7186.538500000001
$ws.Range("M77").Value = -2818.538500000001
$ws.Range("H97").Value = 1939.7693
$ws.Range("I97").Value = 1022.5
$ws.Range("J97").Value = 3407.4
$ws.Range("K97").Value = 1022.5
$ws.Range("L97").Value = 3407.4
$ws.Range("M97").Value = -526.5
$ws.Range("N97").Value = -4399.4
$ws.Range("H132").Value = 2650.08
$ws.Range("I132").Value = 2650.08
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7950.24
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5420.24
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 2334.4443
$ws.Range("I136").Value = 2708.7334
$ws.Range("J136").Value = 463
$ws.Range("K136").Value = 8126.2002
$ws.Range("L136").Value = 1389
$ws.Range("M136").Value = -5576.2002
$ws.Range("N136").Value = -6489
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2006.9166
$ws.Range("I20").Value = 921.8889
$ws.Range("J20").Value = 5262
$ws.Range("K20").Value = 921.8889
$ws.Range("L20").Value = 5262
$ws.Range("M20").Value = -674.8889
$ws.Range("N20").Value = -5756
$ws.Range("H94").Value = 2706046.2
$ws.Range("I94").Value = 3226481.5
$ws.Range("J94").Value = 17131.334
$ws.Range("K94").Value = 3226481.5
$ws.Range("L94").Value = 17131.334
$ws.Range("M94").Value = -3226030.5
$ws.Range("N94").Value = -18033.334
$ws.Range("H134").Value = 4640.639
$ws.Range("I134").Value = 3855.6382
$ws.Range("J134").Value = 7276
$ws.Range("K134").Value = 11566.9146
$ws.Range("L134").Value = 21828
$ws.Range("M134").Value = -9031.9146
$ws.Range("N134").Value = -26898

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3046.2766
$ws.Range("I31").Value = 2073.963
$ws.Range("J31").Value = 4358.9
$ws.Range("K31").Value = 2073.963
$ws.Range("L31").Value = 4358.9
$ws.Range("M31").Value = -1778.963
$ws.Range("N31").Value = -4948.9
$ws.Range("H34").Value = 3046.2766
$ws.Range("I34").Value = 2073.963
$ws.Range("J34").Value = 4358.9
$ws.Range("K34").Value = 2073.963
$ws.Range("L34").Value = 4358.9
$ws.Range("M34").Value = -1871.963
$ws.Range("N34").Value = -4762.9
$ws.Range("H41").Value = 34620.625
$ws.Range("J41").Value = 38995.715
$ws.Range("L41").Value = 38995.715
$ws.Range("N41").Value = -39851.715
$ws.Range("H50").Value = 41999.8
$ws.Range("J50").Value = 41999.8
$ws.Range("L50").Value = 41999.8
$ws.Range("N50").Value = -43249.8
$ws.Range("H58").Value = 1250.5555
$ws.Range("I58").Value = 830.4666999999999
$ws.Range("J58").Value = 2090.7334
$ws.Range("K58").Value = 830.4666999999999
$ws.Range("L58").Value = 2090.7334
$ws.Range("M58").Value = -627.4666999999999
$ws.Range("N58").Value = -2496.7334
$ws.Range("H59").Value = 61535.23
$ws.Range("J59").Value = 61535.23
$ws.Range("L59").Value = 61535.23
$ws.Range("N59").Value = -63825.23
$ws.Range("H60").Value = 23337
$ws.Range("J60").Value = 19985.8
$ws.Range("L60").Value = 19985.8
$ws.Range("N60").Value = -21007.8
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
$ws.Range("H132").Value = 2822.0527
$ws.Range("I132").Value = 2537.68
$ws.Range("J132").Value = 4853.2856
$ws.Range("K132").Value = 7613.039999999999
$ws.Range("L132").Value = 14559.8568
$ws.Range("M132").Value = -5083.039999999999
$ws.Range("N132").Value = -19619.8568
$ws.Range("H134").Value = 3030.0852
$ws.Range("I134").Value = 2381.9302
$ws.Range("J134").Value = 9997.75
$ws.Range("K134").Value = 7145.790599999999
$ws.Range("L134").Value = 29993.25
$ws.Range("M134").Value = -4610.790599999999
$ws.Range("N134").Value = -35063.25
$ws.Range("H136").Value = 1250.5555
$ws.Range("I136").Value = 830.4666999999999
$ws.Range("J136").Value = 2090.7334
$ws.Range("K136").Value = 2491.4001
$ws.Range("L136").Value = 6272.2002
$ws.Range("M136").Value = 58.59990000000016
$ws.Range("N136").Value = -11372.2002

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2421.6
$ws.Range("I80").Value = 2302.4
$ws.Range("K80").Value = 2302.4
$ws.Range("M80").Value = -1304.4
$ws.Range("H83").Value = 2421.6
$ws.Range("I83").Value = 2302.4
$ws.Range("K83").Value = 11512
$ws.Range("M83").Value = -6520
$ws.Range("H132").Value = 2914.5925
$ws.Range("I132").Value = 2919.3333
$ws.Range("J132").Value = 2898
$ws.Range("K132").Value = 8757.999899999999
$ws.Range("L132").Value = 8694
$ws.Range("M132").Value = -6227.999899999999
$ws.Range("N132").Value = -13754

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2404.2896
$ws.Range("I132").Value = 2067.0476
$ws.Range("J132").Value = 2820.8823
$ws.Range("K132").Value = 6201.1428
$ws.Range("L132").Value = 8462.6469
$ws.Range("M132").Value = -3671.1428
$ws.Range("N132").Value = -13522.6469
$ws.Range("H136").Value = 1915.1957
$ws.Range("I136").Value = 1514.4828
$ws.Range("J136").Value = 2598.7646
$ws.Range("K136").Value = 4543.4484
$ws.Range("L136").Value = 7796.293799999999
$ws.Range("M136").Value = -1993.4484
$ws.Range("N136").Value = -12896.2938

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1714.0769
$ws.Range("I126").Value = 1579.1666
$ws.Range("J126").Value = 3333
$ws.Range("K126").Value = 4737.4998
$ws.Range("L126").Value = 9999
$ws.Range("M126").Value = -2267.4998
$ws.Range("N126").Value = -14939
$ws.Range("H132").Value = 2593.043
$ws.Range("I132").Value = 2595.8645
$ws.Range("J132").Value = 2577.9092
$ws.Range("K132").Value = 7787.593500000001
$ws.Range("L132").Value = 7733.7276
$ws.Range("M132").Value = -5257.593500000001
$ws.Range("N132").Value = -12793.7276
$ws.Range("H136").Value = 2430.8333
$ws.Range("I136").Value = 2721.2856
$ws.Range("J136").Value = 1753.1111
$ws.Range("K136").Value = 8163.8568
$ws.Range("L136").Value = 5259.3333
$ws.Range("M136").Value = -5613.8568
$ws.Range("N136").Value = -10359.3333
